# Update cryptos list values (Price and Volume(1h) columns) per upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.446.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +5.17%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.635.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +5.01%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'593.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.58%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'194.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.646"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.35%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.629.94"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +4.99%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.01%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +5.04%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.672"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +3.73%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'58.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +4.47%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +4.98%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'9.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +6.01%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.215.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +4.70%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +6.06%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.635.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +4.69%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'70.502.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +5.11%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +5.13%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +2.16%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +4.77%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'489.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.95%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'19.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +13.25%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'5.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.52%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'4.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.09%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'91.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.23%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +7.70%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'11.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +5.57%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'9.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +5.48%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'33.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +5.48%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +10.70%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +8.82%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'629.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +4.77%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'12.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'65.68"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +3.29%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'41.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +12.75%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.415"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +7.88%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0₃0827"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +9.31%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.94%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.16%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +1.68%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.309.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.69%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +15.42%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +9.63%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0454"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +5.71%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +6.54%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +2.39%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.40%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'9.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +5.60%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'3.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.30%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.07%  "
$ws.Range("E51").Style = "Normal"
